$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns stay formatted as text so that
# values like "6.90" or "60.903.88" are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.903.88"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.917.70"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "590.28"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "146.42"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("D9").Value = "6.90"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "33.56"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "3.400.91"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "60.814.61"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "6.69"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "2.915.10"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "430.21"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "13.36"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").Value = "0.678"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "7.05"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").Value = "81.40"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").Value = "10.95"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "11.85"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +5.02%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").Value = "26.60"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").Value = "0.108"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "0.0₃0858"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").Value = "3.01"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").Value = "40.06"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").Value = "380.59"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").Value = "2.692.10"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "133.42"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "23.79"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -3.62%  "
$ws.Range("E51").Value = "  -0.19%  "
